$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# Add new row 37 data
$ws.Cells.Item(37, 1).Value = 20210526
$ws.Cells.Item(37, 2).Value = 2234.4699999999998
$ws.Cells.Item(37, 3).Value = 2224.4699999999998
$ws.Cells.Item(37, 4).FormulaR1C1 = "=100*(RC[-2]-RC[-1])/RC[-1]"
$ws.Cells.Item(37, 5).Value = 180
$ws.Cells.Item(37, 6).Value = "CRM opened 20210418"

# Update selection / scroll position to match new data extent
$ws.Range("F38").Select()

# Adjust window position (workbookView xWindow/yWindow)
$excel.Windows.Item(1).Left = 1365
$excel.Windows.Item(1).Top = 1710
